# Commit: "added colors for the tasks"
# - Operators sheet: new columns G ("Not evening") and H ("Not task") listing
#   task ids that each operator cannot work.
# - Tasks sheet: new column I ("color") with a hex color per task, plus a
#   couple of tweaked "probability" values.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Operators")
$ws2 = $wb.Worksheets.Item("Tasks")

# --- Sheet1 (Operators): add "Not evening" / "Not task" columns (G, H) ---
$ws1.Cells.Item(1,7).Value = "Not evening"
$ws1.Cells.Item(1,8).Value = "Not task"
$ws1.Cells.Item(2,8).Value = "26,27"
$ws1.Cells.Item(3,8).Value = "8,12,13,26,27"
$ws1.Cells.Item(4,7).Value = "2,9,18,23,30"
$ws1.Cells.Item(4,8).Value = "5,6,19,20"
$ws1.Cells.Item(5,7).Value = "4,9,16"
$ws1.Cells.Item(5,8).Value = "5,6,26,27"
$ws1.Cells.Item(6,7).Value = "9,25"
$ws1.Cells.Item(6,8).Value = "26,27"
$ws1.Cells.Item(7,7).Value = "9,18,23,30"
$ws1.Cells.Item(7,8).Value = "5,6,19,20"
$ws1.Cells.Item(8,7).Value = "2,9,15,23"
$ws1.Cells.Item(8,8).Value = "26,27"
$ws1.Cells.Item(9,7).Value = "2,9"
$ws1.Cells.Item(9,8).Value = "5,6,12,13"
$ws1.Cells.Item(10,7).Value = "4,9,15,22,29"
$ws1.Cells.Item(10,8).Value = "5,6,26,27"
$ws1.Cells.Item(11,7).Value = "2,11,17,23"
$ws1.Cells.Item(11,8).Value = "1,3,12,13,24"
$ws1.Cells.Item(12,7).Value = 23
$ws1.Cells.Item(12,8).Value = "12,13,26,27"
$ws1.Cells.Item(13,7).Value = "11,14,23"
$ws1.Cells.Item(13,8).Value = "5,6,26,27"
$ws1.Cells.Item(14,7).Value = 23
$ws1.Cells.Item(14,8).Value = "26,27,28"
$ws1.Cells.Item(15,8).Value = "11,12,13"
$ws1.Cells.Item(16,7).Value = "14,25,31"
$ws1.Cells.Item(16,8).Value = "9,12,13,26,27"
$ws1.Cells.Item(17,8).Value = "2,3,4,5,6,8,10,11,17,23,25,26,27"
$ws1.Cells.Item(18,7).Value = 8
$ws1.Cells.Item(18,8).Value = "9,10,11,14,15,16,17,18,19,20,21,23,25,26,27,28"
$ws1.Cells.Item(19,7).Value = "2,9"
$ws1.Cells.Item(19,8).Value = "11,12,13,14,26,27"
$ws1.Cells.Item(20,7).Value = "3,23,29"
$ws1.Cells.Item(20,8).Value = "7,17,26,27"
$ws1.Cells.Item(22,7).Value = 10
$ws1.Cells.Item(22,8).Value = "9,17,26,27"
$ws1.Cells.Item(23,7).Value = "10,18,23"
$ws1.Cells.Item(23,8).Value = "5,6,19,20"
$ws1.Cells.Item(24,8).Value = "2,3,5,6,10,21,26,27"
$ws1.Cells.Item(25,7).Value = "17,25"
$ws1.Cells.Item(25,8).Value = "3,10,21,23,26,27,28"
$ws1.Cells.Item(26,8).Value = "8,10,11,12,13,25,26,27"
$ws1.Cells.Item(27,8).Value = "19,20,26,27,30"
$ws1.Cells.Item(28,7).Value = 10
$ws1.Cells.Item(28,8).Value = "19,20,26,27"
$ws1.Cells.Item(29,7).Value = 7
$ws1.Cells.Item(29,8).Value = "3,5,6,17,25,26,27"
$ws1.Cells.Item(30,7).Value = "7,15,23"
$ws1.Cells.Item(30,8).Value = "10,11,26,27,28"
$ws1.Cells.Item(31,7).Value = 7
$ws1.Cells.Item(31,8).Value = "10,11,18,19,20,21,26,27"
$ws1.Cells.Item(32,7).Value = "7,23"
$ws1.Cells.Item(32,8).Value = "10,11,26,27,28"
$ws1.Cells.Item(33,8).Value = "1,26,27"
$ws1.Cells.Item(34,7).Value = 7
$ws1.Cells.Item(34,8).Value = "3,14,17"
$ws1.Cells.Item(35,8).Value = "1,3,4,5,6,9,16,23,26,27"
$ws1.Cells.Item(38,7).Value = 25
$ws1.Cells.Item(38,8).Value = "9,19,20,26,27,28"
$ws1.Cells.Item(40,8).Value = "26,27,28"
$ws1.Cells.Item(41,7).Value = 18
$ws1.Cells.Item(42,8).Value = "1,11"
$ws1.Cells.Item(43,8).Value = "2,9,16,23"

# --- Sheet2 (Tasks): add "color" column (I), update probability values ---
$ws2.Cells.Item(1,9).Value = "color"
$ws2.Cells.Item(2,9).Value = "#F28D9F"
$ws2.Cells.Item(3,9).Value = "#F24452"
$ws2.Cells.Item(4,9).Value = "#F2CB05"
$ws2.Cells.Item(5,9).Value = "#17BF60"
$ws2.Cells.Item(6,9).Value = "#B8D9C4"
$ws2.Cells.Item(7,9).Value = "#F2E205"
$ws2.Cells.Item(8,9).Value = "#F2E205"
$ws2.Cells.Item(9,9).Value = "#F24171"
$ws2.Cells.Item(10,9).Value = "#EAF205"
$ws2.Cells.Item(11,9).Value = "#0BD9D9"
$ws2.Cells.Item(12,9).Value = "#762CBF"
$ws2.Cells.Item(13,9).Value = "#7216F2"
$ws2.Cells.Item(14,9).Value = "#F2EEAC"
$ws2.Cells.Item(15,9).Value = "#EDC4F2"
$ws2.Cells.Item(16,9).Value = "#29A7D9"
$ws2.Cells.Item(17,9).Value = "#F249A6"
$ws2.Cells.Item(18,9).Value = "#50F205"
$ws2.Cells.Item(19,9).Value = "#D9B29C"
$ws2.Cells.Item(14,7).Value = 0.2
$ws2.Cells.Item(19,7).Value = 0.4

# --- Selection / active sheet ---
$ws2.Range("J16").Select()
$ws1.Activate()
$ws1.Range("D1").Select()
